# Apply the LOQ4201 curriculum update.
# The worksheet mirrors each changed value in both column B ("current") and
# column C ("modified, shown in red") since both point at the same shared
# string in the original file. We update both columns to keep them in sync
# with the new shared-string text, matching the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ativação: 01/01/2014 -> 01/01/2021
# Prefix with an apostrophe so Excel stores it as literal text instead of
# auto-converting the date-like string into a numeric date serial (the
# original cell held plain text, not a real date value). Assigning a
# leading-quote string switches on the cell's "quote prefix" flag, which
# would otherwise fork a new cell style; re-pasting the untouched sibling
# cell's formatting over it afterwards restores the original style so only
# the text itself changes, matching the source diff.
$ws.Range("B8").Value = "'01/01/2021"
$ws.Range("C8").Value = "'01/01/2021"
$ws.Range("B7").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("C7").Copy()
$ws.Range("C8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Docentes responsáveis: Marco Antonio Carvalho Pereira -> Herlandí de Souza Andrade
$ws.Range("B13").Value = "11079086 - Herlandí de Souza Andrade"
$ws.Range("C13").Value = "11079086 - Herlandí de Souza Andrade"

# Método:
$ws.Range("B19").Value = "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras"
$ws.Range("C19").Value = "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras"

# Critério:
$ws.Range("B20").Value = "Média Aritmética dos Projetos, Trabalhos e Exercícios realizados no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude) desenvolvidas."
$ws.Range("C20").Value = "Média Aritmética dos Projetos, Trabalhos e Exercícios realizados no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude) desenvolvidas."

# Norma de recuperação:
$ws.Range("B21").Value = "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação."
$ws.Range("C21").Value = "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação."
